# Rename the header row: "<name>_old" -> "<name>_FV2310" and "<name>_new" -> "<name>_FV2404"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value()
    if ($v -like "*_old") {
        $cell.Value = ($v -replace "_old$", "_FV2310")
    } elseif ($v -like "*_new") {
        $cell.Value = ($v -replace "_new$", "_FV2404")
    }
}

# Turn the used range into an Excel Table ("Table1") with an autofilter on the header row
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U63"), [System.Type]::Missing, 1)
$tbl.Name = "Table1"

# Freeze the header row (split below row 1, frozen)
[void]$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
